$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move column B (B1:B52) to column A (A1:A52), carrying values + styles
# across; column C is left untouched. Cut leaves the source range fully
# empty (values + formats), matching the post-edit XML where column B no
# longer has any cell records.
$src = $ws.Range("B1:B52")
$dst = $ws.Range("A1:A52")
$src.Cut($dst)
$ws.Range("B1:B52").Clear()

# Update the active selection to C8, matching the post-edit state.
$ws.Range("C8").Select()
